$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data block (rows 1046-1182) holds weekly Primera/Segunda price records for
# "Coliflor" at "Terminal La Palmera de La Serena". A new week of data needs to be
# added at the top of the block. The new week's prices are carried forward from the
# most recent existing week (the row pair with the latest date, currently at rows
# 1129-1130, dated 45147), with the date advanced by 7 days (one week) to 45154.
# Inserting the two new rows pushes all the existing data rows down by two, which
# matches the observed shift in the workbook.

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R")

$srcRow1 = 1129
$srcRow2 = 1130

# Capture the values of the most recent week (Primera + Segunda) before any rows move.
$vals1 = @{}
$vals2 = @{}
foreach ($c in $cols) {
    $vals1[$c] = $ws.Range($c + $srcRow1).Value2()
    $vals2[$c] = $ws.Range($c + $srcRow2).Value2()
}

# Insert two blank rows at the top of the data block, shifting everything down.
$ws.Rows("1046:1047").Insert()

$newRow1 = 1046
$newRow2 = 1047

foreach ($c in $cols) {
    if ($c -eq "D") {
        $ws.Range($c + $newRow1).Value = $vals1[$c] + 7
        $ws.Range($c + $newRow2).Value = $vals2[$c] + 7
    } else {
        $ws.Range($c + $newRow1).Value = $vals1[$c]
        $ws.Range($c + $newRow2).Value = $vals2[$c]
    }
}
